$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 26 - this shifts the existing data block
# (rows 26-92) down to rows 27-93, preserving all values/formatting.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new week's record.
$ws.Cells.Item(26, 1).Value = 8
$ws.Cells.Item(26, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 44498
$ws.Cells.Item(26, 5).Value = 4
$ws.Cells.Item(26, 6).Value = 100112001
$ws.Cells.Item(26, 7).Value = "Berenjena"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 560
$ws.Cells.Item(26, 11).Value = 8000
$ws.Cells.Item(26, 12).Value = 9000
$ws.Cells.Item(26, 13).Value = 8500
$ws.Cells.Item(26, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(26, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(26, 16).Value = 142
$ws.Cells.Item(26, 17).Value = 60
$ws.Cells.Item(26, 18).Value = "Hortaliza"
